$d = $word.ActiveDocument

# The document currently starts with a single (empty) paragraph that just
# holds the _GoBack bookmark. The edit inserts a brand-new paragraph of
# weekly-report narrative text *before* that paragraph.
$firstPara = $d.Paragraphs.First
$firstRange = $firstPara.Range
$firstRange.InsertParagraphBefore()

# Grab the freshly created (still empty) leading paragraph and populate it
# via a raw WordprocessingML fragment so we get the exact run/proofErr
# layout (the "git" run wrapped in spellStart/spellEnd proofErr markers)
# that the target document contains.
$newPara = $d.Paragraphs.First
$newRange = $newPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t>Weekly report for week ending 2/10/2019. Created 11 stories for Azure DevOps in class and created the timecard and weekly report documents. On Friday after work, an attempt was made to clone the project from home but kept receiving an authentication error, despite using the correct user name/password. Tried deleting credentials in Credential Manager for azure and retry, but no success. On Sunday during the team meeting, for some reason we could only clone them on the computers in the classroom and not in any other computer lab. Much of the meeting time was spent trying to clone the project.</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> Using the command prompt with </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>git</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> commands did not resolve the problem either.</w:t></w:r>' +
  '</w:p>'

[void]$newRange.InsertXML($xml)
